$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix row 2: was missing D2/E2 and had a stray numeric value in B2.
# Now it must match the other rows: A2 keeps its e-mail, B2 gets a new
# password test value, and D2/E2 get the RESULT header-column values
# (copy formatting + content from row 3's D/E cells, the row-2 pattern).
$ws.Range("B2").Value = "Trying_123_.-z"
$ws.Range("D3:E3").Copy()
$ws.Range("D2:E2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D2").Value = "RESULT"

# Update the password test data in column B for the remaining rows.
# (Order matters for how new values are appended to the shared string
# table, so match the original authoring order: B3, B4, B7, then B5.)
$ws.Range("B3").Value = "TRYING"
$ws.Range("B4").Value = "ABC123DFG456"
$ws.Range("B7").Value = "İNvALİD"
$ws.Range("B5").Value = ">>£##£312<<<<3"

# B3 and B4 previously carried extra styling (a custom number format on B3,
# a plain bordered style on B4); both now use the worksheet's default style
# like the rest of column B.
$ws.Range("B3").Style = "Normal"
$ws.Range("B4").Style = "Normal"
